$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before current row 6 (shifts old row6->8, old row7->9)
$ws.Rows.Item(6).Resize(2).EntireRow.Insert() | Out-Null

# Append the new bottom row (row 10): Floor / Z / -22 / =D10*N1
$ws.Range("B10").Value = "Floor"
$ws.Range("C10").Value = "Z"
$ws.Range("D10").Value = -22
$ws.Range("E10").Formula = "=D10*N1"

# New row 7: External East Door Way / Z / 15 / =D7*N1
$ws.Range("B7").Value = "External East Door Way"
$ws.Range("C7").Value = "Z"
$ws.Range("D7").Value = 15
$ws.Range("E7").Formula = "=D7*N1"

# New row 6: External Main Walls / Z / 30 / =D6*N1
$ws.Range("B6").Value = "External Main Walls"
$ws.Range("C6").Value = "Z"
$ws.Range("D6").Value = 30
$ws.Range("E6").Formula = "=D6*N1"
